$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the row 2 student record (order chosen to match the shared-
#     string table layout produced by the original authoring app) ---

# No. (01 -> 02)
$ws.Range("J2").Value = "02"

# E-mail value + its hyperlink
$ws.Range("K2").Value = "thanhb2005691@student.ctu.edu.vn"
$ws.Range("K2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:thanhb2005691@student.ctu.edu.vn")

# Username
$ws.Range("B2").Value = "cam"
# Full name
$ws.Range("D2").Value = "Trương Thị Cam Cam"
# Gender (Nam -> Nữ)
$ws.Range("F2").Value = "Nữ"
# Address (VL -> Vĩnh Long)
$ws.Range("G2").Value = "Vĩnh Long"

# Registration date (serial date, 2024-04-21)
$ws.Range("H2").Value = 45403

# --- Column widths for the newly-sized columns B and C ---
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 9

# --- Selection moved to G2 ---
$ws.Range("G2").Select()
